# Sprint86V01x.xlsx - "Add files via upload" edit
# Mark several Approved test cases as Rejected (with ReasonToReject = "Nil"),
# resize a few columns, and update the sheet view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Approved/Rejected (col I) + ReasonToReject (col J) for the
#     rows that were changed from "Approved" to "Rejected" ---
$rowsToReject = @(2, 4, 5, 7, 8, 10, 12, 33, 49, 56)
foreach ($r in $rowsToReject) {
    $ws.Range("I$r").Value = "Rejected"
    $ws.Range("J$r").Value = "Nil"
}

# --- Resize columns D, E and G ---
$ws.Columns.Item(4).ColumnWidth = 10.0
$ws.Columns.Item(5).ColumnWidth = 9.333333333333332
$ws.Columns.Item(7).ColumnWidth = 13.666666666666668

# --- Update sheet view: zoom to 80% and change the active selection ---
$excel.ActiveWindow.Zoom = 80
$ws.Range("I56").Select()
